$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Table structure changes -------------------------------------------------
# Rename column G header from "Other" to "Note"
$ws.Range("G1").Value = "Note"

# Add a new 8th table column (becomes column H) and name it "Range"
$newcol = $tbl.ListColumns.Add()
$ws.Range("H1").Value = "Range"

# Materialise the (mostly empty) data cells for the new column with the
# same base style ("vertical center / wrap text") used by every other
# data cell in the table (rows 2-18).
$ws.Range("H2:H18").WrapText = $true
$ws.Range("H2:H18").VerticalAlignment = -4108

# --- Column widths -------------------------------------------------------
# xlsx <col width> ends up ~0.8333 higher than the ColumnWidth we set, so
# subtract that offset to land exactly on the target stored widths (61, 13).
$ws.Columns.Item(7).ColumnWidth = 60.166666666666664
$ws.Columns.Item(8).ColumnWidth = 12.166666666666666

# --- Row 2-4: values unchanged, nothing to edit (shared-string reindex only) --

# --- Row 5: password -------------------------------------------------------
$ws.Range("A5").Value = "password"
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 16
$ws.Range("F5").Value = "^[A-Za-z0-9!@#$%^&*]*$"
# Leading "'" = Excel's literal-text quote prefix (text starts with "-", which
# would otherwise look like the start of a formula/number) -- this both keeps
# the text literal and flips on the QuotePrefix style, matching the source.
$ws.Range("G5").Value = "'- Password can contain alphanumeric`n- Password can contains special characters : !@#$%^&*"
$ws.Range("G5").WrapText = $true
$ws.Range("G5").VerticalAlignment = -4108
$ws.Rows.Item(5).RowHeight = 30

# --- Row 6: nickname --------------------------------------------------------
$ws.Range("A6").Value = "nickname"
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 64
$ws.Range("F6").Value = "^[A-Za-z0-9]*$"
$ws.Range("G6").Value = "'- Name can contain alphanumeric"
$ws.Range("G6").WrapText = $true
$ws.Range("G6").VerticalAlignment = -4108

# --- Row 7: created ----------------------------------------------------------
$ws.Range("A7").Value = "created"
$ws.Range("D7").Value = -25200000

# --- Row 8: lastModified -----------------------------------------------------
$ws.Range("A8").Value = "lastModified"
$ws.Range("D8").Value = -25200000

# --- Row 9: joined ------------------------------------------------------------
$ws.Range("A9").Value = "joined"
$ws.Range("D9").Value = -25200000

# --- Row 10: post[title] --------------------------------------------------------
$ws.Range("A10").Value = "post[title]"
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 32
$ws.Range("F10").Value = "^[A-Za-z0-9]*$"

# --- Row 11: post[body] ----------------------------------------------------------
$ws.Range("A11").Value = "post[body]"
$ws.Range("B11").Value = 16
$ws.Range("C11").Value = 512

# --- Row 12: comment[content] ----------------------------------------------------
$ws.Range("A12").Value = "comment[content]"
$ws.Range("B12").Value = 16
$ws.Range("C12").Value = 512

# --- Row 13: account[status] ------------------------------------------------------
$ws.Range("A13").Value = "account[status]"
$ws.Range("H13").Value = "'0: Disabled`n1: Pending`n2: Active`n"
$ws.Range("H13").WrapText = $true
$ws.Range("H13").VerticalAlignment = -4108
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: connection[index] -----------------------------------------------------
$ws.Range("A14").Value = "connection[index]"
$ws.Range("B14").Value = 36
$ws.Range("C14").Value = 36

# --- Selection ----------------------------------------------------------------
$ws.Range("A15").Select()
